# Use default image for officer profiles
# Sets several organizers' imageName (column I) to the shared default
# officer image, and normalizes Willie Chalmers III's filename.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organizers")

# Rows whose profile photo is replaced with the generic default image.
$defaultImageRows = @(3, 14, 20, 22, 23, 24)
foreach ($r in $defaultImageRows) {
    $ws.Cells.Item($r, 9).Value = "officer_default.svg"
}

# Willie Chalmers III: normalize image filename (drop stray space, fix casing).
$ws.Cells.Item(16, 9).Value = "willie_chalmers.jpg"

# Move the active selection to reflect where the editor ended up.
$ws.Activate()
$ws.Range("Q14").Select()
